# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.407.45"
$ws.Range("E2").Value = "  -1.80%  "
$ws.Range("D3").Value = "3.841.63"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'602.45"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "'169.07"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").Value = "3.841.70"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").Value = "'0.167"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  +5.84%  "
$ws.Range("D14").Value = "'37.11"
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").Value = "4.485.84"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "3.833.92"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").Value = "68.450.87"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "'18.48"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").Value = "'7.41"
$ws.Range("E19").Value = "  -2.99%  "
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "'11.08"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").Value = "'470.61"
$ws.Range("D23").Value = "'0.734"
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").Value = "'83.39"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("E26").Value = "  -3.70%  "
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("D28").Value = "'10.08"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "3.991.52"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("D33").Value = "'31.60"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("E35").Value = "  -2.71%  "
$ws.Range("D36").Value = "3.806.86"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("E38").Value = "  +9.13%  "
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  -4.20%  "
$ws.Range("E44").Value = "  -4.40%  "
$ws.Range("D45").Value = "'8.72"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D47").Value = "'415.43"
$ws.Range("E47").Value = "  -4.89%  "
$ws.Range("E49").Value = "  +4.54%  "
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").Value = "'141.51"
$ws.Range("E51").Value = "  -1.50%  "
